# S1_LDATS.docx edit script
# Applies the changes described by the commit's unified diff:
#  1. Merge the "minimal null hypothesis model" sentence into one run (no text change).
#  2. Merge the "if the minimal (species-means) model..." sentence into one run (no text change).
#  3. Insert " (Figure 1)" after "occurring in the mid-1990s".
#  4. Insert "Tables 1-2, " before the first "R2 = .34".
#  5. Merge "no community transition events are detected, as t" into one run (no text change).
#  6. Insert "not shown; " before the second "R2 = ." (the "R2 = . .31 and .35..." one).
#  7. Merge the Simonis et al. LDATS reference into a single run (no text change).

$d = $word.ActiveDocument

# --- 1 ---------------------------------------------------------------
$old1 = "In addition to the models fit using LDA and timeseries analysis, we included as a candidate model a minimal null hypothesis model, in which the prediction for each species in all timesteps is simply the mean proportional abundance over time for that species. "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- 2 ---------------------------------------------------------------
$old2 = "if the minimal (species-means) model was within 2 standard errors of the best-fitting model, "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# --- 3 ---------------------------------------------------------------
$old3 = "occurring in the mid-1990s"
$new3 = "occurring in the mid-1990s (Figure 1)"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# --- 4 ---------------------------------------------------------------
$old4 = " (R2 = .34"
$new4 = " (Tables 1-2, R2 = .34"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# --- 5 ---------------------------------------------------------------
$old5 = "no community transition events are detected, as t"
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $old5, 2) | Out-Null

# --- 6 ---------------------------------------------------------------
$old6 = " (R2 = ..31 and .35 for controls and exclosures, respectively"
$new6 = " (not shown; R2 = ..31 and .35 for controls and exclosures, respectively"
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# --- 7 ---------------------------------------------------------------
$old7 = "Simonis, J. L., E. M. Christensen, D. J. Harris, R. M. Diaz, H. Ye, E. P. White, and S. K. M. Ernest. 2020. LDATS: Latent Dirichlet Allocation Coupled with Time Series Analyses. https://CRAN.R-project.org/package=LDATS"
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $old7, 2) | Out-Null

Write-Output "done"
